$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..J values are identical for the two new rows (4th and 5th
# configuration repeats of the "72500/397" run).
$rowValues = @(72500, 397, 1, 0, 1, 0.0025188916876574307, 0, 0.041666666666666664, 0.9974811083123426)

foreach ($r in 6, 7) {
    $ws.Cells.Item($r, 1).Value = $r - 2

    $col = 2
    foreach ($v in $rowValues) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }

    # Columns K..AB (11..28) hold the per-run stats for runs 1..9, which
    # this configuration never reached -> blank/empty text cells.
    for ($col = 11; $col -le 28; $col++) {
        $ws.Cells.Item($r, $col).Formula = '=""'
    }

    $ws.Cells.Item($r, 29).Value = 0
}

Write-Output "added rows 6 and 7"
